$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 on the "Rules" sheet changes from the text "R40" to the text "1".
# A plain Value/Formula assignment of a numeric-looking string like "1" gets
# auto-coerced to a Number by Excel's normal type inference, but the target
# cell must stay a text cell (t="s"), so we enter it as a text formula and
# then collapse the formula down to its literal text result via a
# copy/paste-special (values only). This preserves the cell's existing
# style (s="23") and leaves a plain text value behind - no formula, no
# number coercion.
$target = $ws.Range("B11")
$target.Formula = "=""1"""
$target.Copy() | Out-Null
$target.PasteSpecial(-4163) | Out-Null      # xlPasteValues
$excel.CutCopyMode = 0
